$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column D width: 28 -> 32 characters (ColumnWidth units need a small offset
# to land exactly on 32 in the stored column width; 31.15 reliably rounds to 32)
$ws.Columns.Item(4).ColumnWidth = 31.15

# Remove all existing hyperlinks first so we can cleanly rebuild them without
# leaving stale relationships behind (per-cell Hyperlinks.Delete() clears the
# whole sheet collection in this engine, so do it once up-front).
$ws.Cells.Hyperlinks.Delete()

# Row 2
$ws.Range("A2").Value2 = "2026-01-07 18:29:46"
$ws.Range("B2").Value2 = "産業機械向けAI異常検知・状態推定システムの開発・導入支援エンジニア募集(AI/エッジ・組み込み)"
$ws.Range("C2").Value2 = "システム開発"
$ws.Range("D2").Value2 = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E2").Value2 = "期限情報なし"
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5450864", "", "", "https://www.lancers.jp/work/detail/5450864")
$ws.Range("F2").Style = "Hyperlink"
$ws.Range("G2").Value2 = 383
$ws.Range("H2").Value2 = "🔥AI,Ai ◆開発"

# Row 3
$ws.Range("A3").Value2 = "2026-01-07 18:29:46"
$ws.Range("B3").Value2 = "【フルタイム】最先端AI(LLM)開発エンジニア募集!新規プロダクトの核となる開発パートナーを募集"
$ws.Range("C3").Value2 = "システム開発"
$ws.Range("D3").Value2 = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Range("E3").Value2 = "期限情報なし"
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5460294", "", "", "https://www.lancers.jp/work/detail/5460294")
$ws.Range("F3").Style = "Hyperlink"
$ws.Range("G3").Value2 = 375
$ws.Range("H3").Value2 = "🔥AI,Ai ◆開発"

# Row 4
$ws.Range("A4").Value2 = "2026-01-07 18:29:46"
$ws.Range("B4").Value2 = "【週5日】法人向け生成AIサービス(RAG・議事録機能)のコア開発を担うリードエンジニア募集"
$ws.Range("C4").Value2 = "システム開発"
$ws.Range("D4").Value2 = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E4").Value2 = "期限情報なし"
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5460267", "", "", "https://www.lancers.jp/work/detail/5460267")
$ws.Range("F4").Style = "Hyperlink"
$ws.Range("G4").Value2 = 375
$ws.Range("H4").Value2 = "🔥AI,Ai ◆開発"

# Row 5
$ws.Range("A5").Value2 = "2026-01-07 18:29:46"
$ws.Range("B5").Value2 = "海外仕入れ元サイト→ツールを動かす為のCSVファイルに週1で自動抽出の制作(自動/スクレイピング)"
$ws.Range("C5").Value2 = "システム開発"
$ws.Range("D5").Value2 = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("E5").Value2 = "期限情報なし"
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5466794", "", "", "https://www.lancers.jp/work/detail/5466794")
$ws.Range("F5").Style = "Hyperlink"
$ws.Range("G5").Value2 = 135
$ws.Range("H5").Value2 = "◆ツール,スクレイピング ◇サイト"

# Row 6
$ws.Range("A6").Value2 = "2026-01-07 18:29:46"
$ws.Range("B6").Value2 = "【急募】宿泊業向けクチコミ対策SaaSのMVP開発"
$ws.Range("C6").Value2 = "システム開発"
$ws.Range("D6").Value2 = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E6").Value2 = "期限情報なし"
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5466852", "", "", "https://www.lancers.jp/work/detail/5466852")
$ws.Range("F6").Style = "Hyperlink"
$ws.Range("G6").Value2 = 75
$ws.Range("H6").Value2 = "◆開発"

# Row 7
$ws.Range("A7").Value2 = "2026-01-07 18:29:46"
$ws.Range("B7").Value2 = "FileMaker開発"
$ws.Range("C7").Value2 = "システム開発"
$ws.Range("D7").Value2 = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E7").Value2 = "期限情報なし"
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5466845", "", "", "https://www.lancers.jp/work/detail/5466845")
$ws.Range("F7").Style = "Hyperlink"
$ws.Range("G7").Value2 = 68
$ws.Range("H7").Value2 = "◆開発"

# Row 8
$ws.Range("A8").Value2 = "2026-01-07 18:29:46"
$ws.Range("B8").Value2 = "Raspberry Piでの開発"
$ws.Range("C8").Value2 = "システム開発"
$ws.Range("D8").Value2 = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E8").Value2 = "期限情報なし"
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5466611", "", "", "https://www.lancers.jp/work/detail/5466611")
$ws.Range("F8").Style = "Hyperlink"
$ws.Range("G8").Value2 = 68
$ws.Range("H8").Value2 = "◆開発"

# Row 9
$ws.Range("A9").Value2 = "2026-01-07 18:29:46"
$ws.Range("B9").Value2 = "【急募】GBP一括投稿システムのインスタ連携改修依頼"
$ws.Range("C9").Value2 = "システム開発"
$ws.Range("D9").Value2 = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E9").Value2 = "期限情報なし"
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5466476", "", "", "https://www.lancers.jp/work/detail/5466476")
$ws.Range("F9").Style = "Hyperlink"
$ws.Range("G9").Value2 = 33
$ws.Range("H9").ClearContents()

# Row 10
$ws.Range("A10").Value2 = "2026-01-07 18:29:46"
$ws.Range("B10").Value2 = "【長期/業務委託】UX改善をリードできるフロント寄り Laravel エンジニア募集(リモート可)"
$ws.Range("C10").Value2 = "システム開発"
$ws.Range("D10").Value2 = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E10").Value2 = "期限情報なし"
$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.lancers.jp/work/detail/5466459", "", "", "https://www.lancers.jp/work/detail/5466459")
$ws.Range("F10").Style = "Hyperlink"
$ws.Range("G10").Value2 = 25
$ws.Range("H10").ClearContents()

# Row 11
$ws.Range("A11").Value2 = "2026-01-07 18:29:46"
$ws.Range("B11").Value2 = "【EC-CUBE】定期購入機能の調査・改修依頼"
$ws.Range("C11").Value2 = "システム開発"
$ws.Range("D11").Value2 = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E11").Value2 = "期限情報なし"
$ws.Hyperlinks.Add($ws.Range("F11"), "https://www.lancers.jp/work/detail/5466925", "", "", "https://www.lancers.jp/work/detail/5466925")
$ws.Range("F11").Style = "Hyperlink"
$ws.Range("G11").Value2 = 18
$ws.Range("H11").ClearContents()

# Row 12
$ws.Range("A12").Value2 = "2026-01-07 18:29:46"
$ws.Range("B12").Value2 = "電気点火装置の回路図作成依頼"
$ws.Range("C12").Value2 = "システム開発"
$ws.Range("D12").Value2 = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E12").Value2 = "期限情報なし"
$ws.Hyperlinks.Add($ws.Range("F12"), "https://www.lancers.jp/work/detail/5466994", "", "", "https://www.lancers.jp/work/detail/5466994")
$ws.Range("F12").Style = "Hyperlink"
$ws.Range("G12").Value2 = 13
$ws.Range("H12").ClearContents()

# Row 13
$ws.Range("A13").Value2 = "2026-01-07 18:29:46"
$ws.Range("B13").Value2 = "【急募】Microsoftドメイン認証&DNS設定のプロを探しています"
$ws.Range("C13").Value2 = "システム開発"
$ws.Range("D13").Value2 = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("E13").Value2 = "期限情報なし"
$ws.Hyperlinks.Add($ws.Range("F13"), "https://www.lancers.jp/work/detail/5466917", "", "", "https://www.lancers.jp/work/detail/5466917")
$ws.Range("F13").Style = "Hyperlink"
$ws.Range("G13").Value2 = 10
$ws.Range("H13").ClearContents()
